$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Manutenção")

# --- Update existing "analisaCandidatos" tracking lines ---
$ws2.Range("C13").Value = "#280 analisaCandidatos"
$ws2.Range("C14").Value = "#530 analisaCandidatos"

# --- New tracking lines appended further down the sheet ---
# (values set in this order so the shared-string table is built in the
#  same sequence as the target workbook)
$ws2.Range("C21").Value = "#530 - iViewNumsDisps!!.setOnTouchListener"
$ws2.Range("C23").Value = "#580 - btnCand.setOnClickListener"
$ws2.Range("C19").Value = "#968 - mostraNumsIguais"
$ws2.Range("C16").Value = "2.3- #1963 analisaCandidatos"
$ws2.Range("C18").Value = "#270 - mostraCelAJogar"
$ws2.Range("C17").Value = "#217 - iViewSudokuBoard!!.setOnTouchListener"

# C18/C19 reuse the same "indent" formatting already used on C13/C14
$ws2.Range("C13").Copy() | Out-Null
$ws2.Range("C18").PasteSpecial(-4122) | Out-Null
$ws2.Range("C19").PasteSpecial(-4122) | Out-Null

# --- Update the saved selection on the "Manutenção" sheet ---
$ws2.Activate() | Out-Null
$ws2.Range("E15").Select() | Out-Null
